$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Locate the "TreatmentRespTab" row (column A = TabName, column B = TabQuery
# holding the raw SQL). Avoids hard-coding the row number.
$targetRow = 0
for ($r = 1; $r -le 9; $r++) {
    $tabName = $ws.Cells.Item($r, 1).Value()
    if ($tabName -eq "TreatmentRespTab") {
        $targetRow = $r
        break
    }
}

$cell = $ws.Cells.Item($targetRow, 2)

# Add the missing "AND trt.treatment_id IS NOT NULL" predicate to the
# Treatment Response query's WHERE clause.
$old = $cell.Value()
$new = $old -replace "std\.dbgap_accession = 'phs002371' AND prt\.sex_at_birth = 'Male'`nORDER BY", "std.dbgap_accession = 'phs002371' AND prt.sex_at_birth = 'Male' AND trt.treatment_id IS NOT NULL`nORDER BY"
$cell.Value = $new

# Leave the cursor/selection on the row that was just edited (matches the
# author's saved cursor position after editing the Treatment Response query).
$ws.Cells.Item($targetRow, 3).Select() | Out-Null
